$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2166
$ws1.Range("G3").Value = 70
$ws1.Range("F4").Value = 887
$ws1.Range("F5").Value = 1462
$ws1.Range("F6").Value = 376

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2166
$ws4.Range("G3").Value = 70
$ws4.Range("F4").Value = 1
$ws4.Range("F6").Value = 887
$ws4.Range("F7").Value = 1462
$ws4.Range("F8").Value = 376
